$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MV data refresh: append the daily observations for 08-09-2021 .. 15-09-2021
# (rows 252-259) to the bottom of the existing series table (A:G).
$newRows = @(
    @("08-09-2021", 2178, 2718, 15718, 2744, 4751, 9127),
    @("09-09-2021", 2152, 2685, 15527, 2710, 4693, 9016),
    @("10-09-2021", 2141, 2672, 15453, 2697, 4671, 8973),
    @("11-09-2021", 2141, 2672, 15453, 2697, 4671, 8973),
    @("12-09-2021", 2141, 2672, 15453, 2697, 4671, 8973),
    @("13-09-2021", 2145, 2676, 15479, 2702, 4679, 8988),
    @("14-09-2021", 2161, 2696, 15591, 2721, 4712, 9053),
    @("15-09-2021", 2163, 2699, 15611, 2725, 4719, 9065)
)

$startRow = 252

# Pre-format the new date cells as text so Excel stores the literal
# "dd-mm-yyyy" strings (like the rest of column A) instead of silently
# re-interpreting them as date serials.
$lastRow = $startRow + $newRows.Count - 1
$dateCol = $ws.Range("A$startRow`:A$lastRow")
$dateCol.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $values[0]
    $ws.Cells.Item($r, 2).Value = $values[1]
    $ws.Cells.Item($r, 3).Value = $values[2]
    $ws.Cells.Item($r, 4).Value = $values[3]
    $ws.Cells.Item($r, 5).Value = $values[4]
    $ws.Cells.Item($r, 6).Value = $values[5]
    $ws.Cells.Item($r, 7).Value = $values[6]
}

# Restore the plain "Normal" style on the date column now that the text
# values are locked in, so the new cells end up unstyled like A2:A251.
$dateCol.Style = "Normal"
